$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that wraps the title text in
#    paragraph 1 ("Searching For Planets In Binary Stars ").
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2. Remove paragraph justification (w:jc val="both") from the three body
#    paragraphs that currently have it, by resetting their alignment back
#    to the (default) left alignment, which drops the <w:jc> element.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Format.Alignment -eq 3) {
        $p.Format.Alignment = 0
    }
}

# 3. Split "Our sample consists of 1100 binary stars..." so the bookmark
#    moves in between "11" and "00", re-inserting the "_GoBack" bookmark
#    there as a zero-length (collapsed) bookmark.
$range = $d.Content
$found = $range.Find.Execute("1100 binary stars", $true, $false, $false,
                              $false, $false, $true, 1, $false, "", 0)
$pt = $d.Range($range.Start + 2, $range.Start + 2)
$d.Bookmarks.Add("_GoBack", $pt)
